# Refresh cryptos list: update price/volume figures and restore text cell types.
# (Row 37/38 coin entries -- OKB and PEPE -- have also swapped places.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'57.748.65"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "'  -3.48%  "
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = "'2.924.35"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "'  -1.62%  "
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').Value = "'  -0.13%  "
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = "'548.68"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "'  -2.99%  "
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = "'130.04"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = "'  +5.17%  "
$ws.Range('E6').Style = "Normal"
$ws.Range('E7').Value = "'  -0.10%  "
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = "'0.509"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = "'  +1.75%  "
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = "'2.918.54"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "'  -1.53%  "
$ws.Range('E9').Style = "Normal"
$ws.Range('E10').Value = "'  -2.33%  "
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = "'4.73"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = "'  -3.79%  "
$ws.Range('E11').Style = "Normal"
$ws.Range('E12').Value = "'  +1.23%  "
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = "'0.0000219"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "'  +0.32%  "
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = "'32.60"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "'  +1.40%  "
$ws.Range('E14').Style = "Normal"
$ws.Range('E15').Value = "'  +1.93%  "
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = "'3.402.76"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "'  -1.96%  "
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = "'6.86"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "'  +5.80%  "
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = "'2.920.21"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Value = "'57.671.51"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = "'  -3.69%  "
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = "'415.60"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = "'  -2.27%  "
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = "'13.26"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = "'  +1.43%  "
$ws.Range('E21').Style = "Normal"
$ws.Range('E22').Value = "'  +3.79%  "
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = "'13.30"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = "'  +4.27%  "
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = "'6.95"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = "'  +0.00%  "
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = "'79.41"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "'  +0.57%  "
$ws.Range('E25').Style = "Normal"
$ws.Range('E27').Value = "'  -0.05%  "
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = "'2.45"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = "'  -2.08%  "
$ws.Range('E28').Style = "Normal"
$ws.Range('E29').Value = "'  +3.83%  "
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = "'7.34"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = "'  +3.07%  "
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = "'25.16"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = "'  +0.69%  "
$ws.Range('E31').Style = "Normal"
$ws.Range('E32').Value = "'  -1.83%  "
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = "'0.0965"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "'  +0.58%  "
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = "'5.67"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "'  +2.42%  "
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = "'0.930"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "'  +1.49%  "
$ws.Range('E35').Style = "Normal"
$ws.Range('E36').Value = "'  +4.07%  "
$ws.Range('E36').Style = "Normal"
$ws.Range('B37').Value = "'PEPE"
$ws.Range('B37').Style = "Normal"
$ws.Range('C37').Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range('C37').Style = "Normal"
$ws.Range('D37').Value = "'0.0₃0688"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = "'  +7.34%  "
$ws.Range('E37').Style = "Normal"
$ws.Range('B38').Value = "'OKB"
$ws.Range('B38').Style = "Normal"
$ws.Range('C38').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('C38').Style = "Normal"
$ws.Range('D38').Value = "'48.16"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "'  -3.90%  "
$ws.Range('E38').Style = "Normal"
$ws.Range('E39').Value = "'  +3.60%  "
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = "'2.58"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "'  +7.65%  "
$ws.Range('E40').Style = "Normal"
$ws.Range('E41').Value = "'  +0.62%  "
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = "'2.701.55"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = "'  +1.79%  "
$ws.Range('E42').Style = "Normal"
$ws.Range('E43').Value = "'  -2.33%  "
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = "'371.37"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = "'  +2.46%  "
$ws.Range('E44').Style = "Normal"
$ws.Range('D46').Value = "'123.73"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "'  +3.05%  "
$ws.Range('E46').Style = "Normal"
$ws.Range('E47').Value = "'  +1.63%  "
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').Value = "'  +0.61%  "
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = "'1.94"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = "'  -0.80%  "
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = "'22.76"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "'  -1.34%  "
$ws.Range('E50').Style = "Normal"
$ws.Range('E51').Value = "'  -0.16%  "
$ws.Range('E51').Style = "Normal"
